# "worked on lectures 2&3 and assignment 2"
#
# Content fix on the "Decision/Chance/Terminal node" legend box (slide 4,
# shape "Content Placeholder 2"): the middle line's text had a typo,
# "Change node", that should read "Chance node".
#
# We locate the paragraph by its current text (robust to index drift) and
# rewrite only its single run's Text so every other run/paragraph
# (formatting, empty spacer paragraphs, etc.) is left completely untouched.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item("Content Placeholder 2")
$tr = $sh.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $text = $para.Text.TrimEnd("`r")
    if ($text -eq "Change node") {
        $para.Runs(1).Text = "Chance node"
    }
}
